$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 70 values (C70, F70)
$ws.Cells.Item(70, 3).Value = 9.830627617743
$ws.Cells.Item(70, 6).Value = 19.830627617743

# Append new rows 74-81
$newRows = @(
    @("2025-02-19", "abs_activity", 10, 9.790047876790231, 0, 19.79004787679023),
    @("2025-02-19", "rel_activity", 0, 5.035750035750036, 0, 5.035750035750036),
    @("2025-02-19", "abs_sleep", 0.2666666666666657, 0, 0, 0.2666666666666657),
    @("2025-02-19", "rel_sleep", 0, 0, 0, 0),
    @("2025-02-20", "abs_activity", 0, 0, 0, 0),
    @("2025-02-20", "rel_activity", 0, 0, 0, 0),
    @("2025-02-20", "abs_sleep", 0, 0, 0, 0),
    @("2025-02-20", "rel_sleep", 0, 0, 0, 0)
)

$startRow = 74
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $row[0]
    $dateCell.Style = "Normal"
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
